$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trim the repo list down to 5 rows: drop the bottom five entries
# (casatablebrowser, casatestdata-large, casafeather, casampi, casa-asap).
$ws.Rows("7:11").Delete()

# Update the repo_name_to_import column (C) for the remaining rows to the
# refreshed repo list.
$ws.Cells.Item(2, 3).Value = "app-n-pak"
$ws.Cells.Item(3, 3).Value = "casa-build-utils"
$ws.Cells.Item(4, 3).Value = "casaplotserver"
$ws.Cells.Item(5, 3).Value = "casashell"
$ws.Cells.Item(6, 3).Value = "casafeather"

# Re-fit the workspace_id / repo_name_to_import columns to their new
# (shorter) content.
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null

# Leave the cursor where it ended up after the edits.
$ws.Range("F11").Select() | Out-Null
